$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing Name values (column A) and Attendance values (column B)
# before overwriting them, since we are shifting Name -> column C and
# inserting a new Roll No column in A.
$names = @()
for ($r = 2; $r -le 7; $r++) {
    $names += $ws.Cells.Item($r, 1).Value2
}

$attendance = @()
for ($r = 2; $r -le 7; $r++) {
    $attendance += $ws.Cells.Item($r, 2).Value2
}

# Header row (set B1 before A1 so the shared-strings table registers
# "attendance" ahead of "Roll No", matching the expected write order)
$ws.Cells.Item(1, 2).Value = "attendance"
$ws.Cells.Item(1, 1).Value = "Roll No"
$ws.Cells.Item(1, 3).Value = "Name"

# Data rows: Roll No (numeric), attendance, Name
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $attendance[$i]
    $ws.Cells.Item($r, 3).Value = $names[$i]
}

$ws.Range("D7").Select()
